# Rename the five worksheets to their new, meaningful names while
# preserving their order, sheetId and relationship ids.
$wb = $excel.ActiveWorkbook

$names = @("Employees", "Managers", "Mentors", "UsefulLinks", "PracticeFormular")

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $wb.Worksheets.Item($i).Name = $names[$i - 1]
}

# The active/first-displayed tab stays on the 5th sheet (now "PracticeFormular"),
# matching the workbook's existing view state.
$wb.Worksheets.Item(5).Activate()
